# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell formatting (bold, border, centered) onto the
# three new header cells, then set their text.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-37: team season record for this roster (66 wins, 48 losses, 1 tie)
for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 29).Value = 66   # AC - Wins
    $ws.Cells.Item($row, 30).Value = 48   # AD - Losses
    $ws.Cells.Item($row, 31).Value = 1    # AE - Ties
}
